$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull / recalculated dSF (column F) values for several rows
$ws.Range("F2").Value = -2
$ws.Range("F3").Value = 7
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = -5
$ws.Range("F8").Value = -3
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = -7
$ws.Range("F14").Value = -4
$ws.Range("F17").Value = -9
$ws.Range("F22").Value = -1
$ws.Range("F23").Value = -7
